# Add three new countries (FRANCE, ITALY, UNITED STATES) to the check_counts
# sheet, keeping the existing alphabetical ordering of the country list and
# pushing all following rows down.
#
# Insertions (by row number in the CURRENT/shifting sheet, processed from
# bottom to top so earlier insert points are unaffected by later ones):
#   - before "URUGUAY"       -> new row: UNITED STATES, 202
#   - before "JAMAICA"       -> new row: ITALY,          153
#   - before "FRENCH GUIANA" -> new row: FRANCE,          205

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert UNITED STATES right before UNITED KINGDOM's successor (URUGUAY), i.e. row 171
$ws.Rows(171).Insert()
$ws.Cells.Item(171, 1).Value2 = "UNITED STATES"
$ws.Cells.Item(171, 2).Value2 = 202

# Insert ITALY right before JAMAICA, i.e. row 83
$ws.Rows(83).Insert()
$ws.Cells.Item(83, 1).Value2 = "ITALY"
$ws.Cells.Item(83, 2).Value2 = 153

# Insert FRANCE right before FRENCH GUIANA, i.e. row 59
$ws.Rows(59).Insert()
$ws.Cells.Item(59, 1).Value2 = "FRANCE"
$ws.Cells.Item(59, 2).Value2 = 205
